$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 368.16666
$ws.Range("I12").Value = 315
$ws.Range("K12").Value = 315
$ws.Range("M12").Value = -145
$ws.Range("H62").Value = 5948.75
$ws.Range("I62").Value = 5948.75
$ws.Range("K62").Value = 5948.75
$ws.Range("M62").Value = -5324.75
$ws.Range("H65").Value = 5948.75
$ws.Range("I65").Value = 5948.75
$ws.Range("K65").Value = 29743.75
$ws.Range("M65").Value = -26623.75
$ws.Range("H76").Value = 6449.5
$ws.Range("I76").Value = 6449.5
$ws.Range("K76").Value = 6449.5
$ws.Range("M76").Value = -6134.5
$ws.Range("H79").Value = 6449.5
$ws.Range("I79").Value = 6449.5
$ws.Range("K79").Value = 6449.5
$ws.Range("M79").Value = -5357.5
$ws.Range("H87").Value = 89999
$ws.Range("J87").Value = 89999
$ws.Range("L87").Value = 89999
$ws.Range("N87").Value = -92495
$ws.Range("H90").Value = 89999
$ws.Range("J90").Value = 89999
$ws.Range("L90").Value = 269997
$ws.Range("N90").Value = -282477
$ws.Range("H132").Value = 22729682
$ws.Range("I132").Value = 27780434
$ws.Range("K132").Value = 83341302
$ws.Range("M132").Value = -83338772
$ws.Range("H137").Value = 17545268
$ws.Range("I137").Value = 27778800
$ws.Range("K137").Value = 83336400
$ws.Range("M137").Value = -83333850

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2235.1538
$ws.Range("I2").Value = 2006.5555
$ws.Range("K2").Value = 2006.5555
$ws.Range("M2").Value = -1893.5555
$ws.Range("H44").Value = 46444
$ws.Range("J44").Value = 55555
$ws.Range("L44").Value = 55555
$ws.Range("N44").Value = -56531
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H55").Value = 54166.25
$ws.Range("J55").Value = 55555
$ws.Range("L55").Value = 55555
$ws.Range("N55").Value = -56185
$ws.Range("H97").Value = 812.619
$ws.Range("I97").Value = 800.7778
$ws.Range("K97").Value = 800.7778
$ws.Range("M97").Value = -304.7778
$ws.Range("H116").Value = 2235.1538
$ws.Range("I116").Value = 2006.5555
$ws.Range("K116").Value = 2006.5555
$ws.Range("M116").Value = 287.4445000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2235.1538
$ws.Range("I3").Value = 2006.5555
$ws.Range("K3").Value = 2006.5555
$ws.Range("M3").Value = -1892.5555
$ws.Range("H86").Value = 12842751
$ws.Range("I86").Value = 21506.666
$ws.Range("J86").Value = 30326268
$ws.Range("K86").Value = 21506.666
$ws.Range("L86").Value = 30326268
$ws.Range("M86").Value = -20383.666
$ws.Range("N86").Value = -30328514
$ws.Range("H89").Value = 12842751
$ws.Range("I89").Value = 21506.666
$ws.Range("J89").Value = 30326268
$ws.Range("K89").Value = 107533.33
$ws.Range("L89").Value = 151631340
$ws.Range("M89").Value = -101917.33
$ws.Range("N89").Value = -151642572
$ws.Range("H94").Value = 1436.4286
$ws.Range("I94").Value = 1507.5
$ws.Range("K94").Value = 1507.5
$ws.Range("M94").Value = -1056.5
$ws.Range("H107").Value = 805.5714
$ws.Range("I107").Value = 805.5714
$ws.Range("K107").Value = 805.5714
$ws.Range("M107").Value = 1114.4286
$ws.Range("H134").Value = 1930.3334
$ws.Range("I134").Value = 1930.3334
$ws.Range("K134").Value = 5791.0002
$ws.Range("M134").Value = -3256.0002
$ws.Range("H141").Value = 84000
$ws.Range("J141").Value = 84000
$ws.Range("L141").Value = 84000
$ws.Range("N141").Value = -94360

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 272.45456
$ws.Range("I107").Value = 259.7
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 259.7
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1660.3
$ws.Range("N107").Value = -4240
$ws.Range("H132").Value = 8702221
$ws.Range("I132").Value = 11119027
$ws.Range("K132").Value = 33357081
$ws.Range("M132").Value = -33354551
$ws.Range("H141").Value = 351007.5
$ws.Range("J141").Value = 351007.5
$ws.Range("L141").Value = 351007.5
$ws.Range("N141").Value = -361367.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2753327.5
$ws.Range("I4").Value = 31844.445
$ws.Range("J4").Value = 15000001
$ws.Range("K4").Value = 95533.33499999999
$ws.Range("L4").Value = 45000003
$ws.Range("M4").Value = -95421.33499999999
$ws.Range("N4").Value = -45000227
$ws.Range("H11").Value = 1226.25
$ws.Range("I11").Value = 1224.091
$ws.Range("K11").Value = 3672.273
$ws.Range("M11").Value = -3532.273
$ws.Range("H107").Value = 367.91666
$ws.Range("I107").Value = 390.85715
$ws.Range("J107").Value = 335.8
$ws.Range("K107").Value = 1172.57145
$ws.Range("L107").Value = 1007.4
$ws.Range("M107").Value = 747.4285500000001
$ws.Range("N107").Value = -4847.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 10044
$ws.Range("J95").Value = 10044
$ws.Range("L95").Value = 10044
$ws.Range("N95").Value = -15536
$ws.Range("H132").Value = 15153873
$ws.Range("I132").Value = 2420
$ws.Range("J132").Value = 41668916
$ws.Range("K132").Value = 7260
$ws.Range("L132").Value = 125006748
$ws.Range("M132").Value = -4730
$ws.Range("N132").Value = -125011808

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 3933.3333
$ws.Range("I122").Value = 3933.3333
$ws.Range("K122").Value = 11799.9999
$ws.Range("M122").Value = -9349.999899999999
$ws.Range("H132").Value = 17928.5
$ws.Range("I132").Value = 11514.2
$ws.Range("K132").Value = 34542.60000000001
$ws.Range("M132").Value = -32012.60000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6829.8887
$ws.Range("I81").Value = 3768.389
$ws.Range("J81").Value = 12952.889
$ws.Range("K81").Value = 7536.778
$ws.Range("L81").Value = 25905.778
$ws.Range("M81").Value = -6475.778
$ws.Range("N81").Value = -28027.778
$ws.Range("H84").Value = 6829.8887
$ws.Range("I84").Value = 3768.389
$ws.Range("J84").Value = 12952.889
$ws.Range("K84").Value = 37683.89
$ws.Range("L84").Value = 129528.89
$ws.Range("M84").Value = -32379.89
$ws.Range("N84").Value = -140136.89
$ws.Range("H100").Value = 1079.8334
$ws.Range("I100").Value = 879.875
$ws.Range("J100").Value = 1479.75
$ws.Range("K100").Value = 1759.75
$ws.Range("L100").Value = 2959.5
$ws.Range("M100").Value = -1218.75
$ws.Range("N100").Value = -4041.5
$ws.Range("H132").Value = 83335070
$ws.Range("I132").Value = 1790.875
$ws.Range("K132").Value = 5372.625
$ws.Range("M132").Value = -2842.625
$ws.Range("H136").Value = 12369
$ws.Range("I136").Value = 15653.875
$ws.Range("K136").Value = 46961.625
$ws.Range("M136").Value = -44411.625
